# edit fault injection code and context table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fault-injection code / context table text edits -----------------------
# Row 6 ("low glucose") and row 8 ("over glucose") conditions were reworded
# from the generic bg_target comparison to the loaded_glucose fixture values
# used by the fault-injection tests.
$ws.Range("B6").Value = "loaded_glucose < 120"
$ws.Range("B8").Value = "loaded_glucose > 120"

# --- Column width tweaks (context table got re-laid-out / resized) ---------
$ws.Columns.Item(1).ColumnWidth = 7.666666666666667
$ws.Columns.Item(2).ColumnWidth = 27.333333333333332
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 49.166666666666664

# --- Window / view state ----------------------------------------------------
$excel.ActiveWindow.TabRatio = 0.991

# Active cell moved to D9 after the edits.
$ws.Range("D9").Select() | Out-Null
